# Apply the update described by the diff:
#  - Remove the oldest 10 years of data (2000年..2009年), shifting the
#    remaining rows (2010年..2020年) up into rows 2..12.
#  - Append a new row for 2021年 with full data (previously missing).
#  - Append a new row for 2022年 with only the "参加养老保险人数" (column D)
#    value populated; the other data columns are left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for 2000年 through 2009年 (original rows 2-11). This
# shifts 2010年 (was row 12) up to row 2, ... 2020年 (was row 22) up to
# row 12.
$ws.Rows("2:11").Delete()

# --- Row 13: 2021年 -------------------------------------------------
# Clone the formatting of the preceding data row (now row 12, 2020年) so
# the new row keeps the same cell style (border/bold/alignment) without
# introducing a brand new style entry.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 31101.5002
$ws.Range("C13").Value = 11126.4955
$ws.Range("D13").Value = 48074.0377
$ws.Range("E13").Value = 34917.0713
$ws.Range("F13").Value = 13156.9664

# --- Row 14: 2022年 -------------------------------------------------
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)

$ws.Range("A14").Value = "2022年"
$ws.Range("D14").Value = 50349
